# Update cryptocurrency Price (D) and Volume(1h) (E) columns per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "27.860.71"
$c.Style = "Normal"

$c = $ws.Range("E2")
$c.Value = "  -0.89%  "
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.Value = "1.906.68"
$c.Style = "Normal"

$c = $ws.Range("E3")
$c.Value = "  -0.22%  "
$c.Style = "Normal"

$c = $ws.Range("E4")
$c.Value = "  -0.17%  "
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.Value = "'313.41"
$c.Style = "Normal"

$c = $ws.Range("E5")
$c.Value = "  -1.12%  "
$c.Style = "Normal"

$c = $ws.Range("E6")
$c.Value = "  -0.20%  "
$c.Style = "Normal"

$c = $ws.Range("D7")
$c.Value = "'0.5010"
$c.Style = "Normal"

$c = $ws.Range("E7")
$c.Value = "  +4.03%  "
$c.Style = "Normal"

$c = $ws.Range("E8")
$c.Value = "  -0.06%  "
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.Value = "'0.07283"
$c.Style = "Normal"

$c = $ws.Range("E9")
$c.Value = "  -1.10%  "
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.Value = "'0.9075"
$c.Style = "Normal"

$c = $ws.Range("E10")
$c.Value = "  -2.94%  "
$c.Style = "Normal"

$c = $ws.Range("D11")
$c.Value = "'20.87"
$c.Style = "Normal"

$c = $ws.Range("E11")
$c.Value = "  +0.06%  "
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.Value = "1.944.58"
$c.Style = "Normal"

$c = $ws.Range("E12")
$c.Value = "  +1.76%  "
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.Value = "'0.07666"
$c.Style = "Normal"

$c = $ws.Range("E13")
$c.Value = "  -1.63%  "
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.Value = "'5.483"
$c.Style = "Normal"

$c = $ws.Range("E14")
$c.Value = "  -0.73%  "
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.Value = "'91.61"
$c.Style = "Normal"

$c = $ws.Range("E15")
$c.Value = "  -0.13%  "
$c.Style = "Normal"

$c = $ws.Range("E16")
$c.Value = "  -0.25%  "
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.Value = "'0.000008722"
$c.Style = "Normal"

$c = $ws.Range("E17")
$c.Value = "  -1.29%  "
$c.Style = "Normal"

$c = $ws.Range("E18")
$c.Value = "  -0.09%  "
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.Value = "27.900.42"
$c.Style = "Normal"

$c = $ws.Range("E19")
$c.Value = "  -0.86%  "
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.Value = "'5.172"
$c.Style = "Normal"

$c = $ws.Range("E21")
$c.Value = "  -0.02%  "
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.Value = "'10.83"
$c.Style = "Normal"

$c = $ws.Range("E22")
$c.Value = "  -0.85%  "
$c.Style = "Normal"

$c = $ws.Range("D23")
$c.Value = "'6.606"
$c.Style = "Normal"

$c = $ws.Range("E23")
$c.Value = "  -0.71%  "
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.Value = "'154.43"
$c.Style = "Normal"

$c = $ws.Range("D25")
$c.Value = "'1.879"
$c.Style = "Normal"

$c = $ws.Range("E25")
$c.Value = "  -2.27%  "
$c.Style = "Normal"

$c = $ws.Range("E26")
$c.Value = "  +5.29%  "
$c.Style = "Normal"

$c = $ws.Range("E27")
$c.Value = "  -0.89%  "
$c.Style = "Normal"

$c = $ws.Range("D28")
$c.Value = "'115.29"
$c.Style = "Normal"

$c = $ws.Range("E28")
$c.Value = "  -1.21%  "
$c.Style = "Normal"

$c = $ws.Range("D29")
$c.Value = "'4.912"
$c.Style = "Normal"

$c = $ws.Range("E29")
$c.Value = "  -1.05%  "
$c.Style = "Normal"

$c = $ws.Range("D30")
$c.Value = "'0.08975"
$c.Style = "Normal"

$c = $ws.Range("D31")
$c.Value = "'3.207"
$c.Style = "Normal"

$c = $ws.Range("E31")
$c.Value = "  -3.90%  "
$c.Style = "Normal"

$c = $ws.Range("E32")
$c.Value = "  -2.10%  "
$c.Style = "Normal"

$c = $ws.Range("D33")
$c.Value = "'0.7679"
$c.Style = "Normal"

$c = $ws.Range("E33")
$c.Value = "  -1.20%  "
$c.Style = "Normal"

$c = $ws.Range("D34")
$c.Value = "'4.651"
$c.Style = "Normal"

$c = $ws.Range("E34")
$c.Value = "  -0.84%  "
$c.Style = "Normal"

$c = $ws.Range("D35")
$c.Value = "'0.02063"
$c.Style = "Normal"

$c = $ws.Range("E35")
$c.Value = "  -0.06%  "
$c.Style = "Normal"

$c = $ws.Range("D36")
$c.Value = "'2.551"
$c.Style = "Normal"

$c = $ws.Range("E36")
$c.Value = "  -3.75%  "
$c.Style = "Normal"

$c = $ws.Range("D37")
$c.Value = "'0.5587"
$c.Style = "Normal"

$c = $ws.Range("E37")
$c.Value = "  +1.37%  "
$c.Style = "Normal"

$c = $ws.Range("D38")
$c.Value = "'1.094"
$c.Style = "Normal"

$c = $ws.Range("E38")
$c.Value = "  -1.39%  "
$c.Style = "Normal"

$c = $ws.Range("D39")
$c.Value = "'3.016"
$c.Style = "Normal"

$c = $ws.Range("E39")
$c.Value = "  +1.01%  "
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.Value = "'0.05250"
$c.Style = "Normal"

$c = $ws.Range("E40")
$c.Value = "  -1.45%  "
$c.Style = "Normal"

$c = $ws.Range("D41")
$c.Value = "'6.958"
$c.Style = "Normal"

$c = $ws.Range("E41")
$c.Value = "  -0.94%  "
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.Value = "'8.496"
$c.Style = "Normal"

$c = $ws.Range("E42")
$c.Value = "  -0.08%  "
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.Value = "'0.1512"
$c.Style = "Normal"

$c = $ws.Range("E43")
$c.Value = "  -1.19%  "
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.Value = "'111.53"
$c.Style = "Normal"

$c = $ws.Range("E44")
$c.Value = "  +2.94%  "
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.Value = "'0.4830"
$c.Style = "Normal"

$c = $ws.Range("E45")
$c.Value = "  -0.40%  "
$c.Style = "Normal"

$c = $ws.Range("D46")
$c.Value = "'10.53"
$c.Style = "Normal"

$c = $ws.Range("E46")
$c.Value = "  -2.17%  "
$c.Style = "Normal"

$c = $ws.Range("E47")
$c.Value = "  -0.23%  "
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.Value = "'1.631"
$c.Style = "Normal"

$c = $ws.Range("E48")
$c.Value = "  -1.79%  "
$c.Style = "Normal"

$c = $ws.Range("D49")
$c.Value = "'67.56"
$c.Style = "Normal"

$c = $ws.Range("E49")
$c.Value = "  -1.11%  "
$c.Style = "Normal"

$c = $ws.Range("D50")
$c.Value = "'0.06067"
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.Value = "'0.9011"
$c.Style = "Normal"

$c = $ws.Range("E51")
$c.Value = "  -0.21%  "
$c.Style = "Normal"

